# Regenerate the handback-status report timestamps.
# (commit: "Generate Report for Handback")
#
# Updates the four timestamp cells that get refreshed each time the
# handback-status report is (re)generated:
#   - Overview!G2               "Latest HO Xliff Generate Date"
#   - zh-cn!H2                  "Correspond Handoff Datetime"
#   - zh-cn!K2                  "Correspond Handback DateTime"
#   - de-de!H2                  "Correspond Handoff Datetime"
#   - de-de!K2                  "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 21:05:28"

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 21:05:23"
$wsZhCn.Range("K2").Value = "2016-08-17 21:05:43"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 21:05:28"
$wsDeDe.Range("K2").Value = "2016-08-17 21:05:51"
